$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New bit-field mini table in D2:E3 (OpEnI = b2, OpEnQ = b3) ---
# Write in this order so the shared-string table gets indices
# 139=b2, 140=OpEnI, 141=b3, 142=OpEnQ (matches target workbook).
$ws.Range("E2").Value = "b2"
$ws.Range("E3").Value = "OpEnI"
$ws.Range("D2").Value = "b3"
$ws.Range("D3").Value = "OpEnQ"

# --- Consolidate the SPI instr/data register sub-tables (rows 11-16) ---
# Row 12/13 gain a D column (label) alongside the existing E column,
# and E12/E13 switch from SPI_INSTR/SPI指令 to SPI_DATA/SPI数据.
$ws.Range("D12").Value = "SPI_INSTR"
$ws.Range("D13").Value = "SPI指令"
$ws.Range("E12").Value = "SPI_DATA"
$ws.Range("E13").Value = "SPI数据"

# Match the bordered "label" style already used by neighboring cells
# (copy formatting only, so no stray duplicate style gets created).
$ws.Range("F1").Copy() | Out-Null
$ws.Range("D12:D13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Rows 15/16 no longer carry their own E-column label (now merged up
# into rows 12/13), so clear them out entirely.
$ws.Range("E15").Clear()
$ws.Range("E16").Clear()

# --- View state: scroll position + active selection ---
$excel.Goto($ws.Range("A7"), $true) | Out-Null
$ws.Range("I19").Select() | Out-Null
